# Update the handoff/handback datetimes on the generated handback report
# to reflect the latest run ("Generate Report for handback").

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-17 15:17:49"
$wsZhCn.Range("G2").Value = "2016-01-17 15:18:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-17 15:17:59"
$wsDeDe.Range("G2").Value = "2016-01-17 15:18:51"
